# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
# F2: 1237 -> 1239
# F3: 74   -> 75

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1239
    $ws.Range("F3").Value = 75
}
